$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows (5 and 6) of login data.
# The order in which NEW distinct string values are first written controls
# the order they receive in the shared-strings table, so write them in the
# exact sequence needed to reproduce shared string indices 5, 6, 7.
$ws.Range("B5").Value = "Password2"
$ws.Range("A6").Value = "r iti@gmail.com"
$ws.Range("A5").Value = " riti@gmail.com "
$ws.Range("B6").Value = "Password2"

# Give the new rows the same bordered look as the rest of the table by
# copying the format from an existing bordered cell (B1) onto A5:B6.
$ws.Range("B1").Copy()
$ws.Range("A5:B6").PasteSpecial(-4122)

# Column A of the new rows also got its fill explicitly toggled off
# (an explicit "No Fill"), which marks the fill as applied in the
# generated style even though it still renders as no fill.
$ws.Range("A5:A6").Interior.ColorIndex = -4142

# Update the active selection to F7, matching the saved view state.
$null = $ws.Range("F7").Select()
